# Adds a new "Installations" worksheet containing a small summary table
# (Contrator / Count of installations) backed by an Excel Table named
# "Table2", matching the target workbook produced by the commit.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet, placed after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Installations"

# --- Populate the cells ---
# (Values are entered in this particular order so that any newly created
#  shared-string table entries come out in the same order as the target
#  workbook: Contrator, Deezlo, Count of installations, I.P.E, Nimba.)
$ws.Range("A2").Value = "Contrator"
$ws.Range("A3").Value = "Deezlo"
$ws.Range("B2").Value = "Count of installations"
$ws.Range("A4").Value = "I.P.E"
$ws.Range("A5").Value = "Nimba"
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0

# --- Column widths (matches author's 11.25 / 20.125 character widths) ---
$ws.Columns.Item(1).ColumnWidth = 10.5
$ws.Columns.Item(2).ColumnWidth = 19.333333333333332

# --- Turn the range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $ws.Range("A2:B5"),
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table2"
$tbl.TableStyle = "TableStyleLight8"

# --- Match the recorded selection/active cell on the new sheet ---
$ws.Range("K8").Select() | Out-Null

Write-Host "Added 'Installations' worksheet with Table2."
